$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New diary entries for 2/13 and 2/16, added as rows 21 and 22 (previously
# blank template rows). Row 20's formatting is cloned onto row 21 (same
# column style pattern: date / time / participants / goal / achievements /
# reflection / mood), then the values are overwritten.
# ---------------------------------------------------------------------------

# Clone formatting from the last filled-in row (20) onto row 21 so the new
# entry picks up the same style indices (date format, bold/italic labels,
# wrap text, etc.)
$ws.Range("A20:G20").Copy() | Out-Null
$ws.Range("A21:G21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 21 - 2/13/2020 entry (43874 = 2/13/2020 serial date number)
$ws.Range("A21").Value = 43874
$ws.Range("B21").Value = "5:00 -7:50 pm"
$ws.Range("C21").Value = "N/A"
$ws.Range("D21").Value = "big picture and more key expert practices"
$ws.Range("E21").Value = "We learned about various stakeholders important to the development of software and how experts work along different levels of abstrction and how they prioritize work"
$ws.Range("F21").Value = "I really found the key expert practice ""do something else"" very helpful. Often I get stuck on trying to understand how a piece of code works and I will fixate on it. I think I needed to hear that it's good practice to stop and search other areas of the code or to do somethin entirely different while your mind sorts out hte problem. "
$ws.Range("G21").Value = "Good"
$ws.Rows("21").RowHeight = 153

# Row 22 - 2/16/2020 entry (keeps the row's existing format, only the date
# column needs the date-number style that row 21/20 already use)
$ws.Range("A20").Copy() | Out-Null
$ws.Range("A22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 43877 = 2/16/2020 serial date number
$ws.Range("A22").Value = 43877
$ws.Range("B22").Value = "3:00 - 7:30pm"
$ws.Range("C22").Value = "Chris, Jay, Rafael"
$ws.Range("D22").Value = "Worked on finding the stakeholders for latest project"
$ws.Range("E22").Value = "we were able to accomplish our goal of finding stakeholders by searching thorugh forums, github, and documentation"
$ws.Range("F22").Value = "I learned the value of documentation, forums, and github. We were lucky to have such an engaged community with our OS project. By searching through forums we learned that our OS project is actively engaged with its users by implementnig features and reporting/fixing bugs. It was great to see such an active role from the developr side and makes me appreciate how much care goes into someones software"
$ws.Rows("22").RowHeight = 204

# ---------------------------------------------------------------------------
# Selection / view state at time of save
# ---------------------------------------------------------------------------
$ws.Range("F22").Select() | Out-Null
